$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name / link updates (row 51: WhiteBITCoin -> dogwifhat)
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"

# Price column updates - force text storage so values like "1.00"/"18.97"
# are preserved exactly as strings (matching the original inlineStr cells)
# instead of being auto-coerced into numeric values by Excel.
$dCells = @("D2","D3","D5","D6","D10","D13","D14","D15","D16","D17","D19","D21","D22","D23","D24","D25","D27","D29","D31","D32","D34","D35","D38","D39","D43","D44","D46","D48","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D2").Value = "62.857.27"
$ws.Range("D3").Value = "2.579.54"
$ws.Range("D5").Value = "582.20"
$ws.Range("D6").Value = "146.44"
$ws.Range("D10").Value = "5.64"
$ws.Range("D13").Value = "27.18"
$ws.Range("D14").Value = "3.041.88"
$ws.Range("D15").Value = "62.723.20"
$ws.Range("D16").Value = "0.0000147"
$ws.Range("D17").Value = "2.586.49"
$ws.Range("D19").Value = "341.72"
$ws.Range("D21").Value = "6.65"
$ws.Range("D22").Value = "0.999"
$ws.Range("D23").Value = "5.65"
$ws.Range("D24").Value = "67.02"
$ws.Range("D25").Value = "2.711.45"
$ws.Range("D27").Value = "1.59"
$ws.Range("D29").Value = "7.86"
$ws.Range("D31").Value = "1.44"
$ws.Range("D32").Value = "1.92"
$ws.Range("D34").Value = "460.23"
$ws.Range("D35").Value = "175.10"
$ws.Range("D38").Value = "0.398"
$ws.Range("D39").Value = "18.97"
$ws.Range("D43").Value = "157.85"
$ws.Range("D44").Value = "3.75"
$ws.Range("D46").Value = "21.03"
$ws.Range("D48").Value = "0.0964"
$ws.Range("D50").Value = "18.39"
$ws.Range("D51").Value = "1.71"
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# Volume(1h) percentage-change text updates
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +2.09%  "
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("E10").Value = "  +2.43%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("E16").Value = "  +3.01%  "
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("E21").Value = "  -1.40%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("E29").Value = "  +7.12%  "
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("E31").Value = "  -3.22%  "
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("E34").Value = "  +13.23%  "
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("E36").Value = "  +3.56%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("E40").Value = "  +3.98%  "
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("E43").Value = "  +4.35%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  +5.50%  "
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("E51").Value = "  +0.33%  "
